$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 carries the workbook default (unstyled) cell style; used to strip the
# temporary text-number-format style back off each edited Price cell so the
# cell keeps its original (absent) style index after forcing a text value.
$defaultStyle = $ws.Range("A1").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.014.33"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.760.24"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "644.96"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +2.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.67"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.760.04"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  -1.95%  "

$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("E9").Value = "  +0.48%  "

$ws.Range("E10").Value = "  -2.05%  "

$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.88"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +3.46%  "

$ws.Range("E13").Value = "  -5.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.89"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  -3.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.393.10"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -1.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.757.80"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  -0.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.026.70"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.68"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  -2.95%  "

$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("E20").Value = "  -2.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.62"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.56"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -1.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.706"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000143"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -6.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.66"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -2.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.25"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +1.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.22"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +1.25%  "

$ws.Range("E28").Value = "  -3.95%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.908.20"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  -1.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.69"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.25"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.13"
$ws.Range("D33").Style = $defaultStyle

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.51"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -2.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.172"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +15.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.715.10"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.77"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -3.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  -3.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.76"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -2.73%  "

$ws.Range("E41").Value = "  -6.51%  "

$ws.Range("E42").Value = "  +0.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.955"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -2.70%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.93"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +2.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.98"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  +2.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.72"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.11"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  +0.30%  "

$ws.Range("E49").Value = "  -2.29%  "

$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("E51").Value = "  -1.48%  "
